# Backup before MoClo code restructure
# Restructure the Echo liquid-handler transfer list: the four Level-1 DNA
# parts (promoter / RBS / CDS / terminator) get their real BioBrick/iGEM
# part numbers instead of the generic placeholder names, and the transfer
# table is extended from 20 to 30 rows to cover destination wells A5 and A6
# (the two pTU1 backbone assemblies) as well as the remaining combinations
# for destination wells A3 and A6 that use a doubled (500 nL) transfer
# volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full transfer table for rows 2-31: row, UID, Source Plate Name, Source
# Plate Type, Source Well, Destination Plate Name, Destination Plate Type,
# Destination Well, Transfer Volume, Reagent. Rows 2-17 (the four DNA parts
# x wells A1-A4) keep their layout but the Reagent names move from the
# generic placeholders to real BioBrick/iGEM part numbers; rows 18-31 are
# rebuilt/extended to cover destination wells A5 and A6 (incl. the two
# pTU1 backbone assemblies).
$data = @(
    @(2, 1, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "BBa_J23119"),
    @(3, 2, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "BBa_B0031"),
    @(4, 3, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "BBa_K1323010"),
    @(5, 4, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "BBa_B1001"),
    @(6, 5, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 250, "BBa_J23119"),
    @(7, 6, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 250, "BBa_B0031"),
    @(8, 7, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 250, "BBa_K1323010"),
    @(9, 8, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 250, "BBa_B1001"),
    @(10, 9, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 500, "BBa_J23119"),
    @(11, 10, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 500, "BBa_B0031"),
    @(12, 11, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 500, "BBa_K1323010"),
    @(13, 12, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 500, "BBa_B1001"),
    @(14, 13, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_J23119"),
    @(15, 14, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_B0031"),
    @(16, 15, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_K1323010"),
    @(17, 16, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "BBa_B1001"),
    @(18, 17, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A5", 250, "BBa_J23119"),
    @(19, 18, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A5", 250, "BBa_B0031"),
    @(20, 19, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A5", 250, "BBa_K1323010"),
    @(21, 20, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A5", 250, "BBa_B1001"),
    @(22, 21, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A1", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A6", 500, "BBa_J23119"),
    @(23, 22, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A2", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A6", 500, "BBa_B0031"),
    @(24, 23, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A3", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A6", 500, "BBa_K1323010"),
    @(25, 24, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A4", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A6", 500, "BBa_B1001"),
    @(26, 25, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A5", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A1", 250, "pTU1-A-lacZ"),
    @(27, 26, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A5", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A2", 500, "pTU1-A-lacZ"),
    @(28, 27, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A5", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A3", 250, "pTU1-A-lacZ"),
    @(29, 28, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A6", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A4", 250, "pTU1-B-lacZ"),
    @(30, 29, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A6", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A5", 500, "pTU1-B-lacZ"),
    @(31, 30, "level 1 384 source plate (DNA components)", "384LDV_AQ_B", "A6", "384-Well Level 1 MoClo output plate", "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)", "A6", 250, "pTU1-B-lacZ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}
